# Mustafizur Rahman (Rajasthan Royals) batting log — add the "matchNo"
# column and the remaining scraped match rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet title: generic "Sheet1" -> player name.
$ws.Name = "Mustafizur Rahman"

# 2) Make room for the new first column ("matchNo"); this shifts the
#    existing teamName..result data from A:L to B:M.
$ws.Columns("A:A").Insert()

# Helper: write a value as plain text. Values that look like numbers
# ("0", "114.28", ...) get a leading quote so Excel keeps them as text
# (matching the sheet's existing numberStoredAsText convention) instead
# of silently converting them to the Number type.
function Write-Text($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -eq "") {
        # Keep the cell present as an empty *string* (not just blank/cleared)
        # by quote-prefixing nothing -- mirrors the sheet's existing empty
        # "states" cell, which is an empty inline string rather than a
        # truly-blank cell.
        $cell.Value = "'"
        return
    }
    $looksNumeric = $text -match '^-?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

function Write-Row($ws, $row, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = $i + 1
        $val = $values[$i]
        Write-Text $ws $row $col $val
    }
}

# 3) Header row: add the new "matchNo" label in A1 (B1:M1 already hold
#    the original headers after the column insert above).
Write-Text $ws 1 1 "matchNo"

# 4) The single pre-existing data row (now B2:M2, after the column
#    insert) is one of the five scraped matches — move it down to row 6
#    and tag it with its match number.
$ws.Range("B2:M2").Cut($ws.Range("B6:M6"))
Write-Text $ws 6 1 "12th"
# Re-assert the (empty) "states" cell so it stays a proper empty string,
# consistent with the same column on the other rows, rather than the
# untyped blank cell the cut/paste leaves behind.
Write-Text $ws 6 4 ""

# 5) Remaining four matches, newest first, as rows 2-5.
$matches = @(
    @("51st", "Rajasthan Royals", "Mustafizur Rahman", "", "8", "7", "0", "1", "114.28", "Mumbai Indians", "Sharjah", "October 05", "Mumbai won by 8 wickets (with 70 balls remaining)"),
    @("32nd", "Rajasthan Royals", "Mustafizur Rahman", "", "0", "0", "0", "0", "-", "Punjab Kings", "Dubai (DSC)", "September 21", "Royals won by 2 runs"),
    @("16th", "Rajasthan Royals", "Mustafizur Rahman", "", "0", "0", "0", "0", "-", "Royal Challengers Bangalore", "Wankhede", "April 22", "RCB won by 10 wickets (with 21 balls remaining)"),
    @("54th", "Rajasthan Royals", "Mustafizur Rahman", "", "0", "3", "0", "0", "0.00", "Kolkata Knight Riders", "Sharjah", "October 07", "KKR won by 86 runs")
)

$row = 2
foreach ($m in $matches) {
    Write-Row $ws $row $m
    $row++
}
